$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 658 ("「感情／感覚」" entry) - all subsequent rows shift up by one.
$ws.Rows.Item(658).Delete()
